$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# sheet1: insert two new leading columns ("Unnamed: 0", "Unnamed: 0.1") before
# the existing "max" column, and add matching data in row 2.
# ---------------------------------------------------------------------------

# Style the new header cells (C1, D1) like the existing header cell B1:
# bold font, thin box border, centered horizontal, top vertical alignment.
$hdrRange1 = $ws1.Range("C1:D1")
$hdrRange1.Font.Bold = $true
$hdrRange1.HorizontalAlignment = -4108
$hdrRange1.VerticalAlignment = -4160
$hdrRange1.Borders.LineStyle = 1

# Shift header text: B1 becomes "Unnamed: 0", C1 "Unnamed: 0.1", D1 keeps "max".
$ws1.Range("D1").Value = "max"
$ws1.Range("B1").Value = "Unnamed: 0"
$ws1.Range("C1").Value = "Unnamed: 0.1"

# Row 2 data: A2 (index column, keeps its style) becomes 0, B2 stays 0, and
# the new C2/D2 cells hold the shifted "min"/extra values.
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 0.3541453925452973

# ---------------------------------------------------------------------------
# sheet2: insert two new leading columns ("Unnamed: 0", "Unnamed: 0.1")
# before the existing "max"/"min" columns; the old label column (a..k)
# moves from column A into column C, and column A becomes a numeric row
# index.
# ---------------------------------------------------------------------------

$hdrRange2 = $ws2.Range("D1:E1")
$hdrRange2.Font.Bold = $true
$hdrRange2.HorizontalAlignment = -4108
$hdrRange2.VerticalAlignment = -4160
$hdrRange2.Borders.LineStyle = 1

$ws2.Range("D1").Value = "max"
$ws2.Range("E1").Value = "min"
$ws2.Range("B1").Value = "Unnamed: 0"
$ws2.Range("C1").Value = "Unnamed: 0.1"

$labels = @("a", "b", "c", "d", "e", "f", "g", "h", "k")

for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2

    # Capture the pre-edit B/C values (the old "max"/"min" data columns) before
    # they get overwritten by the new column layout.
    $oldB = $ws2.Range("B$r").Value2
    $oldC = $ws2.Range("C$r").Value2

    # New layout: A=row index, B=row index (dup), C=old label (a..k),
    # D/E=old max/min data, shifted two columns to the right.
    $ws2.Range("D$r").Value = $oldB
    $ws2.Range("E$r").Value = $oldC
    $ws2.Range("C$r").Value = $labels[$i]
    $ws2.Range("A$r").Value = $i
    $ws2.Range("B$r").Value = $i
}
